$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27 (pasture_nr): update columns E-H
$ws.Range("E27").Value = 22.440128326416016
$ws.Range("F27").Value = 21.062820434570313
$ws.Range("G27").Value = 22.436120986938477
$ws.Range("H27").Value = 20.461940765380859

# Row 28 (CRP_nr): update columns C-H
$ws.Range("C28").Value = 97.661148071289063
$ws.Range("D28").Value = 81.848159790039063
$ws.Range("E28").Value = 72.973960876464844
$ws.Range("F28").Value = 67.612373352050781
$ws.Range("G28").Value = 62.581783294677734
$ws.Range("H28").Value = 63.415077209472656
